$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate the formatting of the last existing data rows (66 & 67) down onto
# the two new rows (68 & 69) so the new cells pick up the same styles
# (date format on column A, numeric format on C:Z) instead of creating new
# style/numFmt entries.
$ws.Range("A66:Z67").Copy()
$ws.Range("A68:Z69").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$stationA = "四方坪站充电量(kw)"
$stationB = "高岭站充电量(kw)"

$row68 = @(45934, $stationA, 439.15499999999997, 829.96, 795.94500000000005, 519.45400000000006, 211.40900000000002, 976.58500000000015, 276.95000000000005, 289.524, 130.24, 115.77099999999999, 127.62599999999999, 159.93200000000002, 474.80999999999983, 1021.1700000000002, 311.03399999999999, 415.59, 235.69, 202.285, 132.74, 221.97999999999996, 56.41, 118.96000000000001, 0.08, 175.184)

$row69 = @(45934, $stationB, 412.27199999999993, 394.13800000000003, 159.5, 70.581000000000003, 173.84100000000001, 173.45999999999998, 179.70400000000001, 38.918999999999997, 311.149, 151.29, 132.96199999999999, 194.55699999999999, 340.01999999999992, 234.58600000000001, 258.04800000000006, 240.87100000000001, 158.02199999999999, 78.882999999999996, 50.030999999999999, 86.652000000000001, 28.125, 92.718999999999994, 4.5860000000000003, 22.934000000000001)

for ($c = 0; $c -lt $row68.Length; $c++) {
    $ws.Cells.Item(68, $c + 1).Value = $row68[$c]
}

for ($c = 0; $c -lt $row69.Length; $c++) {
    $ws.Cells.Item(69, $c + 1).Value = $row69[$c]
}

$ws.Range("F73").Select()
